$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.25566983192155135
$ws.Range("A2").Value = -0.05574138369895465
$ws.Range("A3").Value = -0.003999999977628121
$ws.Range("A4").Value = -0.007999999959332982
$ws.Range("A5").Value = -0.0029999999764243057
$ws.Range("A6").Value = -0.0019999999739894747
$ws.Range("A7").Value = -0.00999999994161005
$ws.Range("A8").Value = -0.009999999939480642
$ws.Range("A9").Value = -0.001999999968281152
$ws.Range("A10").Value = -0.0019999999658999457
$ws.Range("A11").Value = -0.002999999961720512
$ws.Range("A12").Value = -0.003499999959152511
$ws.Range("A13").Value = -0.0034999999566558415
$ws.Range("A14").Value = -0.007999999938209434
$ws.Range("A15").Value = -0.0009999999652929858
$ws.Range("A16").Value = -0.0019999999609252583
$ws.Range("A17").Value = -0.001999999960276888
$ws.Range("A18").Value = -0.003999999952262634
$ws.Range("A19").Value = -0.0039999999825663934
$ws.Range("A20").Value = -0.003999999978585578
$ws.Range("A21").Value = -0.003999999977765789
$ws.Range("A22").Value = -0.026146649338259387
$ws.Range("A23").Value = -0.004999999972381097
$ws.Range("A24").Value = -0.019999999911152422
$ws.Range("A25").Value = -0.019999999910014665
$ws.Range("A26").Value = -0.0024999999719614863
$ws.Range("A27").Value = -0.002499999971326883
$ws.Range("A28").Value = -0.0019999999704403137
$ws.Range("A29").Value = -0.006999999948833491
$ws.Range("A30").Value = -0.05999999973993653
$ws.Range("A31").Value = -0.006999999950110691
$ws.Range("A32").Value = 0.025781250700084968
$ws.Range("A33").Value = 0.013750536784373324
